$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B4").Value = "SOIC-8_208mil"
$ws.Range("A2:A6").ClearFormats()
$ws.Range("C2:D6").ClearFormats()
$ws.Range("B5").Select()
